# Actualización automática 2025-10-30 08:30:09
# Updates the "PIEDRA SINTERIZADA" sale for PALATE CHUCARALAO JOSE ISRAEL
# and the "PORCELANATO" sale for TAMAYO CONDO LUIS ALFREDO, plus the
# downstream monthly / compliance totals that depend on them.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L24").Value = 5179.53
$ws1.Range("M32").Value = 2536.39
$ws1.Range("L37").Value = "5 de 35"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F24").Value = 5179.53
$ws2.Range("F32").Value = 2536.39
$ws2.Range("F37").Value = 43255.3

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 11: PIEDRA SINTERIZADA
$ws3.Range("D11").Value = 7434.9
$ws3.Range("E11").Value = -4512.67541814726
$ws3.Range("F11").Value = 2.544260302979912

# Row 12: PORCELANATO
$ws3.Range("D12").Value = 33857.19
$ws3.Range("E12").Value = -12155.92
$ws3.Range("F12").Value = 1.560147862314049

# Row 14: TOTAL
$ws3.Range("D14").Value = 44657.2
$ws3.Range("E14").Value = -8071.632762818175
$ws3.Range("F14").Value = 1.220623414432536
